# Auto-generated edit script: updates Leve Profit calculation columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per updated market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 154.375
$ws.Range("I2").Value = 154.375
$ws.Range("K2").Value = 154.375
$ws.Range("M2").Value = -41.375

$ws.Range("H6").Value = 110879.3
$ws.Range("I6").Value = 112088.22
$ws.Range("K6").Value = 336264.66
$ws.Range("M6").Value = -336152.66

$ws.Range("H7").Value = 9124.5
$ws.Range("J7").Value = 9124.5
$ws.Range("L7").Value = 9124.5
$ws.Range("N7").Value = -9348.5

$ws.Range("H14").Value = 9124.5
$ws.Range("J14").Value = 9124.5
$ws.Range("L14").Value = 9124.5
$ws.Range("N14").Value = -9506.5

$ws.Range("H15").Value = 1836.5217
$ws.Range("I15").Value = 1836.5217
$ws.Range("K15").Value = 5509.5651
$ws.Range("M15").Value = -5340.5651

$ws.Range("H17").Value = 418301.1
$ws.Range("I17").Value = 2500
$ws.Range("J17").Value = 436379.4
$ws.Range("K17").Value = 7500
$ws.Range("L17").Value = 1309138.2
$ws.Range("M17").Value = -7332
$ws.Range("N17").Value = -1309474.2

$ws.Range("H28").Value = 1062.5
$ws.Range("I28").Value = 576
$ws.Range("J28").Value = 2278.75
$ws.Range("K28").Value = 576
$ws.Range("L28").Value = 2278.75
$ws.Range("M28").Value = -91
$ws.Range("N28").Value = -3248.75

$ws.Range("H40").Value = 3991
$ws.Range("J40").Value = 1989.6666
$ws.Range("L40").Value = 1989.6666
$ws.Range("N40").Value = -2339.6666

$ws.Range("H43").Value = 54999
$ws.Range("J43").Value = 9999
$ws.Range("L43").Value = 9999
$ws.Range("N43").Value = -10137

$ws.Range("H55").Value = 788.2
$ws.Range("I55").Value = 520.5
$ws.Range("J55").Value = 966.6667
$ws.Range("K55").Value = 520.5
$ws.Range("L55").Value = 966.6667
$ws.Range("M55").Value = -306.5
$ws.Range("N55").Value = -1394.6667

$ws.Range("H74").Value = 4833.3335
$ws.Range("I74").Value = 4478.143
$ws.Range("J74").Value = 5330.6
$ws.Range("K74").Value = 4478.143
$ws.Range("L74").Value = 5330.6
$ws.Range("M74").Value = -3542.143
$ws.Range("N74").Value = -7202.6

$ws.Range("H77").Value = 4833.3335
$ws.Range("I77").Value = 4478.143
$ws.Range("J77").Value = 5330.6
$ws.Range("K77").Value = 22390.715
$ws.Range("L77").Value = 26653
$ws.Range("M77").Value = -17710.715
$ws.Range("N77").Value = -36013

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H95").Value = 12000
$ws.Range("J95").Value = 12000
$ws.Range("L95").Value = 12000
$ws.Range("N95").Value = -17492

$ws.Range("H116").Value = 29249.25
$ws.Range("J116").Value = 5666
$ws.Range("L116").Value = 5666
$ws.Range("N116").Value = -12550

$ws.Range("H132").Value = 2571.0833
$ws.Range("I132").Value = 2884.5
$ws.Range("J132").Value = 1004
$ws.Range("K132").Value = 8653.5
$ws.Range("L132").Value = 3012
$ws.Range("M132").Value = -6123.5
$ws.Range("N132").Value = -8072

$ws.Range("H137").Value = 38593.344
$ws.Range("I137").Value = 56158.7
$ws.Range("J137").Value = 9317.75
$ws.Range("K137").Value = 168476.1
$ws.Range("L137").Value = 27953.25
$ws.Range("M137").Value = -165926.1
$ws.Range("N137").Value = -33053.25

$ws.Range("H138").Value = 2632.87
$ws.Range("I138").Value = 1460.3334
$ws.Range("J138").Value = 2944.557
$ws.Range("K138").Value = 4381.0002
$ws.Range("L138").Value = 8833.670999999998
$ws.Range("M138").Value = 758.9997999999996
$ws.Range("N138").Value = -19113.671

$ws.Range("H141").Value = 843.1053
$ws.Range("I141").Value = 778.8333
$ws.Range("K141").Value = 2336.4999
$ws.Range("M141").Value = 2843.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 427.5
$ws.Range("I5").Value = 203.33333
$ws.Range("K5").Value = 203.33333
$ws.Range("M5").Value = -91.33332999999999

$ws.Range("H61").Value = 3687.1428
$ws.Range("I61").Value = 2086.2778
$ws.Range("J61").Value = 5382.1763
$ws.Range("K61").Value = 2086.2778
$ws.Range("L61").Value = 5382.1763
$ws.Range("M61").Value = -1874.2778
$ws.Range("N61").Value = -5806.1763

$ws.Range("H122").Value = 4356.25
$ws.Range("I122").Value = 3393.7273
$ws.Range("K122").Value = 10181.1819
$ws.Range("M122").Value = -7731.1819

$ws.Range("H132").Value = 1704.5294
$ws.Range("I132").Value = 1704.5294
$ws.Range("K132").Value = 5113.5882
$ws.Range("M132").Value = -2583.5882

$ws.Range("H136").Value = 3687.1428
$ws.Range("I136").Value = 2086.2778
$ws.Range("J136").Value = 5382.1763
$ws.Range("K136").Value = 6258.8334
$ws.Range("L136").Value = 16146.5289
$ws.Range("M136").Value = -3708.8334
$ws.Range("N136").Value = -21246.5289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 427.5
$ws.Range("I4").Value = 203.33333
$ws.Range("K4").Value = 203.33333
$ws.Range("M4").Value = -88.33332999999999

$ws.Range("H14").Value = 99999
$ws.Range("J14").Value = 99999
$ws.Range("L14").Value = 99999
$ws.Range("N14").Value = -100343

$ws.Range("H22").Value = 50389
$ws.Range("I22").Value = 779
$ws.Range("K22").Value = 779
$ws.Range("M22").Value = -606

$ws.Range("H26").Value = 209100.7
$ws.Range("I26").Value = 209100.7
$ws.Range("K26").Value = 209100.7
$ws.Range("M26").Value = -208808.7

$ws.Range("H75").Value = 45107
$ws.Range("I75").Value = 40214
$ws.Range("K75").Value = 40214
$ws.Range("M75").Value = -39278

$ws.Range("H76").Value = 48999
$ws.Range("J76").Value = 48999
$ws.Range("L76").Value = 48999
$ws.Range("N76").Value = -49629

$ws.Range("H78").Value = 45107
$ws.Range("I78").Value = 40214
$ws.Range("K78").Value = 120642
$ws.Range("M78").Value = -115962

$ws.Range("H79").Value = 48999
$ws.Range("J79").Value = 48999
$ws.Range("L79").Value = 48999
$ws.Range("N79").Value = -51183

$ws.Range("H99").Value = 3317.606
$ws.Range("I99").Value = 3160.0715
$ws.Range("K99").Value = 3160.0715
$ws.Range("M99").Value = -1662.0715

$ws.Range("H103").Value = 20639.8
$ws.Range("J103").Value = 20639.8
$ws.Range("L103").Value = 20639.8
$ws.Range("N103").Value = -22983.8

$ws.Range("H134").Value = 3758.4707
$ws.Range("I134").Value = 3330
$ws.Range("K134").Value = 9990
$ws.Range("M134").Value = -7455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 473.89474
$ws.Range("I22").Value = 399.93332
$ws.Range("K22").Value = 399.93332
$ws.Range("M22").Value = -49.93331999999998

$ws.Range("H62").Value = 6200
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 6200
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H68").Value = 31952.111
$ws.Range("J68").Value = 33446.125
$ws.Range("L68").Value = 33446.125
$ws.Range("N68").Value = -34944.125

$ws.Range("H71").Value = 31952.111
$ws.Range("J71").Value = 33446.125
$ws.Range("L71").Value = 100338.375
$ws.Range("N71").Value = -107826.375

$ws.Range("H99").Value = 487590.88
$ws.Range("J99").Value = 44998
$ws.Range("L99").Value = 44998
$ws.Range("N99").Value = -47994

$ws.Range("H107").Value = 4847.2896
$ws.Range("J107").Value = 5142.8286
$ws.Range("L107").Value = 5142.8286
$ws.Range("N107").Value = -8982.8286

$ws.Range("H122").Value = 4995
$ws.Range("I122").Value = 4995
$ws.Range("K122").Value = 14985
$ws.Range("M122").Value = -12535

$ws.Range("H126").Value = 487590.88
$ws.Range("J126").Value = 44998
$ws.Range("L126").Value = 134994
$ws.Range("N126").Value = -139934

$ws.Range("H132").Value = 3772.7144
$ws.Range("I132").Value = 3372.353
$ws.Range("K132").Value = 10117.059
$ws.Range("M132").Value = -7587.059000000001

$ws.Range("H134").Value = 6075.2705
$ws.Range("J134").Value = 2786
$ws.Range("L134").Value = 8358
$ws.Range("N134").Value = -13428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 697.8
$ws.Range("J9").Value = 499.5
$ws.Range("L9").Value = 1498.5
$ws.Range("N9").Value = -1946.5

$ws.Range("H12").Value = 195.84616
$ws.Range("J12").Value = 227.125
$ws.Range("L12").Value = 681.375
$ws.Range("N12").Value = -1027.375

$ws.Range("H36").Value = 781.6667
$ws.Range("I36").Value = 698
$ws.Range("J36").Value = 1200
$ws.Range("K36").Value = 2094
$ws.Range("L36").Value = 3600
$ws.Range("M36").Value = -1925
$ws.Range("N36").Value = -3938

$ws.Range("H46").Value = 6112291
$ws.Range("I46").Value = 802.5714
$ws.Range("K46").Value = 2407.7142
$ws.Range("M46").Value = -2316.7142

$ws.Range("H76").Value = 350001000
$ws.Range("I76").Value = 350001000
$ws.Range("K76").Value = 1050003000
$ws.Range("M76").Value = -1050002617

$ws.Range("H79").Value = 350001000
$ws.Range("I79").Value = 350001000
$ws.Range("K79").Value = 1050003000
$ws.Range("M79").Value = -1050001674

$ws.Range("H80").Value = 6496.909
$ws.Range("J80").Value = 6774.1113
$ws.Range("L80").Value = 20322.3339
$ws.Range("N80").Value = -22194.3339

$ws.Range("H83").Value = 6496.909
$ws.Range("J83").Value = 6774.1113
$ws.Range("L83").Value = 60967.00169999999
$ws.Range("N83").Value = -70327.0017

$ws.Range("H114").Value = 561.3077
$ws.Range("I114").Value = 480.2
$ws.Range("J114").Value = 831.6667
$ws.Range("K114").Value = 1440.6
$ws.Range("L114").Value = 2495.0001
$ws.Range("M114").Value = 1813.4
$ws.Range("N114").Value = -9003.000100000001

$ws.Range("H117").Value = 788.86664
$ws.Range("J117").Value = 938.9
$ws.Range("L117").Value = 2816.7
$ws.Range("N117").Value = -9700.7

$ws.Range("H121").Value = 47621256
$ws.Range("J121").Value = 2734.0715
$ws.Range("L121").Value = 8202.2145
$ws.Range("N121").Value = -10822.2145

$ws.Range("H122").Value = 1856.8462
$ws.Range("J122").Value = 2500.7144
$ws.Range("L122").Value = 22506.4296
$ws.Range("N122").Value = -27406.4296

$ws.Range("H129").Value = 19232328
$ws.Range("I129").Value = 27778586
$ws.Range("J129").Value = 3248.5
$ws.Range("K129").Value = 83335758
$ws.Range("L129").Value = 9745.5
$ws.Range("M129").Value = -83330758
$ws.Range("N129").Value = -19745.5

$ws.Range("H132").Value = 5891.241
$ws.Range("I132").Value = 7864.7896
$ws.Range("J132").Value = 2141.5
$ws.Range("K132").Value = 70783.1064
$ws.Range("L132").Value = 19273.5
$ws.Range("M132").Value = -68253.1064
$ws.Range("N132").Value = -24333.5

$ws.Range("H137").Value = 4707.5
$ws.Range("I137").Value = 4030
$ws.Range("J137").Value = 4933.3335
$ws.Range("K137").Value = 12090
$ws.Range("L137").Value = 14800.0005
$ws.Range("M137").Value = -6990
$ws.Range("N137").Value = -25000.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.86667
$ws.Range("I2").Value = 110
$ws.Range("K2").Value = 110
$ws.Range("M2").Value = 3

$ws.Range("H70").Value = 7921.0347
$ws.Range("I70").Value = 7869.4585
$ws.Range("K70").Value = 7869.4585
$ws.Range("M70").Value = -7599.4585

$ws.Range("H73").Value = 7921.0347
$ws.Range("I73").Value = 7869.4585
$ws.Range("K73").Value = 7869.4585
$ws.Range("M73").Value = -6933.4585

$ws.Range("H126").Value = 9100
$ws.Range("I126").Value = 9797.143
$ws.Range("K126").Value = 29391.429
$ws.Range("M126").Value = -26921.429

$ws.Range("H132").Value = 60485.156
$ws.Range("I132").Value = 70613.25
$ws.Range("J132").Value = 6468.6665
$ws.Range("K132").Value = 211839.75
$ws.Range("L132").Value = 19405.9995
$ws.Range("M132").Value = -209309.75
$ws.Range("N132").Value = -24465.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 55555810
$ws.Range("I16").Value = 55555810
$ws.Range("K16").Value = 55555810
$ws.Range("M16").Value = -55555640

$ws.Range("H40").Value = 20747.428
$ws.Range("I40").Value = 21828.076
$ws.Range("K40").Value = 21828.076
$ws.Range("M40").Value = -21692.076

$ws.Range("H46").Value = 5609.4443
$ws.Range("I46").Value = 4995
$ws.Range("J46").Value = 5916.6665
$ws.Range("K46").Value = 4995
$ws.Range("L46").Value = 5916.6665
$ws.Range("M46").Value = -4807
$ws.Range("N46").Value = -6292.6665

$ws.Range("H68").Value = 7758.643
$ws.Range("I68").Value = 9768.143
$ws.Range("K68").Value = 9768.143
$ws.Range("M68").Value = -9019.143

$ws.Range("H71").Value = 7758.643
$ws.Range("I71").Value = 9768.143
$ws.Range("K71").Value = 48840.715
$ws.Range("M71").Value = -45096.715

$ws.Range("H93").Value = 2847.4348
$ws.Range("I93").Value = 2067.0908
$ws.Range("J93").Value = 3562.75
$ws.Range("K93").Value = 2067.0908
$ws.Range("L93").Value = 3562.75
$ws.Range("M93").Value = -819.0908
$ws.Range("N93").Value = -6058.75

$ws.Range("H136").Value = 1954.6296
$ws.Range("I136").Value = 1843.421
$ws.Range("K136").Value = 5530.263
$ws.Range("M136").Value = -2980.263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 566.5
$ws.Range("I107").Value = 566.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1699.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 220.5
$ws.Range("N107").ClearContents()

$ws.Range("H117").Value = 38332.668
$ws.Range("J117").Value = 38332.668
$ws.Range("L117").Value = 38332.668
$ws.Range("N117").Value = -47510.668

$ws.Range("H126").Value = 3087.7778
$ws.Range("I126").Value = 3162.5
$ws.Range("J126").Value = 2490
$ws.Range("K126").Value = 9487.5
$ws.Range("L126").Value = 7470
$ws.Range("M126").Value = -7017.5
$ws.Range("N126").Value = -12410

$ws.Range("H136").Value = 401222.62
$ws.Range("I136").Value = 455880.4
$ws.Range("J136").Value = 399
$ws.Range("K136").Value = 1367641.2
$ws.Range("L136").Value = 1197
$ws.Range("M136").Value = -1365091.2
$ws.Range("N136").Value = -6297

$ws.Range("H137").Value = 89800
$ws.Range("J137").Value = 89800
$ws.Range("L137").Value = 89800
$ws.Range("N137").Value = -100000
